$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 18.65696907043457
$ws.Range("D2").Value = 196

$ws.Range("C3").Value = 17.5330638885498
$ws.Range("D3").Value = 174

$ws.Range("C4").Value = 17.13275909423828
$ws.Range("D4").Value = 174

$ws.Range("C5").Value = 17.27819442749023
$ws.Range("D5").Value = 193

$ws.Range("C6").Value = 17.37117767333984
$ws.Range("D6").Value = 182
